# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet (price1, price2, discount1, ... boosters),
# prefix each command name in column A (rows 2..N, i.e. everything except
# the "Name" header in row 1) with the worksheet's own name followed by a
# space, e.g. "Step4 Seed" -> "free1 Step4 Seed".

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $prefix = $name + " "

    $usedRange = $ws.UsedRange
    $firstRow = $usedRange.Row
    $lastRow = $firstRow + $usedRange.Rows.Count - 1

    # Row 1 is the header ("Name"); data rows start at row 2.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Text
        if ($current -ne "") {
            $cell.Value = $prefix + $current
        }
    }
}
